# Logged Week 16 and performed season sim from Week 17
#
# 1) Append the new week's play-by-play yardage/return logs to the
#    space-separated strings that track every play of the season.
# 2) Refresh the season-to-date aggregate totals on OFF / DEF / ST /
#    TURNS / PEN to reflect the newly logged week plus the simulated
#    remainder of the season (weeks 17+).

$wb = $excel.ActiveWorkbook

function Append-PlayLog($ws, $cellRef, $newValues) {
    $cell = $ws.Range($cellRef)
    $existing = $cell.Value()
    $cell.Value = ($existing + " " + $newValues).Trim()
}

# ---------------------------------------------------------------------
# YDS sheet: per-play rush/pass yardage logs (offense + defense)
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets.Item("YDS")
Append-PlayLog $ydsWs "B2" "1 5 18 4 -2 2 -1 1 15 1 -1 0 -2 2 8 -3 9 2 3 5 -2 1"
Append-PlayLog $ydsWs "B3" "6 1 9 9 14 25 -1 4 5 9 8 6 9 6 21 34 9 5 19 19 14 6 10 6 30 16"
Append-PlayLog $ydsWs "C2" "8 2 5 7 6 1 10 8 9 0 16 4 2 5 4 4 10 3 11 17 3 4 0 4 3 3 7 -1 1 2 0"
Append-PlayLog $ydsWs "C3" "5 4 19 11 8 6 4 3 6 6 10 4 7 16 11 8 12 7 37 6 7"

# ---------------------------------------------------------------------
# ST sheet: kickoff-return and punt-return per-play logs
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets.Item("ST")
Append-PlayLog $stWs "B6" "9 0 7"
Append-PlayLog $stWs "D3" "37 47 63 39 40"
Append-PlayLog $stWs "D4" "0 14 0 0 61"
Append-PlayLog $stWs "D5" "0 0"

# ---------------------------------------------------------------------
# OFF sheet: season aggregate totals (row 2 = RATT, row 3 = PATT)
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets.Item("OFF")
$offWs.Range("C2").Value = 202
$offWs.Range("D2").Value = 10
$offWs.Range("F2").Value = 66
$offWs.Range("G2").Value = 56
$offWs.Range("I2").Value = 3
$offWs.Range("L2").Value = 250
$offWs.Range("M2").Value = 163
$offWs.Range("O2").Value = 19
$offWs.Range("P2").Value = 11
$offWs.Range("Q2").Value = 455

$offWs.Range("C3").Value = 151
$offWs.Range("E3").Value = 44
$offWs.Range("F3").Value = 79
$offWs.Range("G3").Value = 30
$offWs.Range("H3").Value = 35
$offWs.Range("I3").Value = 64
$offWs.Range("J3").Value = 50
$offWs.Range("N3").Value = 11

# ---------------------------------------------------------------------
# DEF sheet: season aggregate totals (row 2 = RATT, row 3 = PATT)
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets.Item("DEF")
$defWs.Range("C2").Value = 214
$defWs.Range("F2").Value = 59
$defWs.Range("G2").Value = 55
$defWs.Range("I2").Value = 5
$defWs.Range("J2").Value = 23
$defWs.Range("L2").Value = 251
$defWs.Range("M2").Value = 157
$defWs.Range("O2").Value = 23
$defWs.Range("Q2").Value = 456

$defWs.Range("C3").Value = 163
$defWs.Range("F3").Value = 113
$defWs.Range("G3").Value = 39
$defWs.Range("I3").Value = 59
$defWs.Range("J3").Value = 46

# ---------------------------------------------------------------------
# ST sheet: season aggregate totals (row 2 = KO/PT/XP/FG counts, row 3 = Home)
# ---------------------------------------------------------------------
$stWs.Range("B2").Value = 86
$stWs.Range("D2").Value = 65
$stWs.Range("F2").Value = 73
$stWs.Range("G2").Value = 65
$stWs.Range("J2").Value = 68
$stWs.Range("K2").Value = 57
$stWs.Range("L2").Value = 36
$stWs.Range("M2").Value = 23

$stWs.Range("B3").Value = 71

# ---------------------------------------------------------------------
# TURNS sheet: season aggregate totals (row 2 = Home, row 3 = Road)
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets.Item("TURNS")
$turnsWs.Range("B2").Value = 7
$turnsWs.Range("C2").Value = 6
$turnsWs.Range("D2").Value = 7
$turnsWs.Range("E2").Value = 7

$turnsWs.Range("D3").Value = 4
$turnsWs.Range("E3").Value = 6

# ---------------------------------------------------------------------
# PEN sheet: season aggregate totals (row 2 = False start)
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets.Item("PEN")
$penWs.Range("B2").Value = 16
